$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$find = "Administrator, Miss Dina Nasr"
$replace = "Miss Dina Nasr, Administrator"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $find) {
        $cell.Value2 = $replace
    }
}
